$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (F column) counts
$wsExhibitions = $wb.Worksheets.Item("展览")
$wsExhibitions.Range("F3").Value = 12169
$wsExhibitions.Range("F6").Value = 374
$wsExhibitions.Range("F8").Value = 12086
$wsExhibitions.Range("F9").Value = 508
$wsExhibitions.Range("F10").Value = 1188
$wsExhibitions.Range("F11").Value = 117
$wsExhibitions.Range("F14").Value = 5955
$wsExhibitions.Range("F15").Value = 138
$wsExhibitions.Range("F16").Value = 3565
$wsExhibitions.Range("F17").Value = 207
$wsExhibitions.Range("F18").Value = 33

# Sheet "全部类型" (All types) - update "想去人数" (F column) counts
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 12170
$wsAll.Range("F8").Value = 374
$wsAll.Range("F10").Value = 12086
$wsAll.Range("F11").Value = 508
$wsAll.Range("F12").Value = 1188
$wsAll.Range("F13").Value = 117
$wsAll.Range("F17").Value = 5955
$wsAll.Range("F18").Value = 138
$wsAll.Range("F19").Value = 3565
$wsAll.Range("F20").Value = 207
$wsAll.Range("F21").Value = 33
